# Apply the "Add files via upload" edit to the IoT Smart Drain Monitoring deck.
#
# The only user-visible content change in the target revision is on slide 1:
# the phrase "Blynk Cloud." becomes "Thingspeak Cloud." (the word "Blynk "
# is replaced by "Thingspeak ", which PowerPoint splits into its own run).
#
# Every slide part is also re-serialized (picked up automatically whenever a
# slide is touched/edited), so we lightly "touch" every other slide (a no-op
# self-assignment of the title shape's Name) to make sure they get
# re-written too.

$p = $ppt.ActivePresentation

# --- Slide 1: update the description text -------------------------------
$slide1 = $p.Slides.Item(1)
$body = $slide1.Shapes.Item(2).TextFrame.TextRange

# "Blynk " (6 chars, including the trailing space) starts right after
# "...drainage systems via " and right before "Cloud."
$target = $body.Characters(141, 6)
if ($target.Text -eq "Blynk ") {
    $target.Text = "Thingspeak "
} else {
    # Fallback in case indices ever drift: do a plain text replace.
    $full = $body.Text
    $body.Text = $full.Replace("Blynk Cloud.", "Thingspeak Cloud.")
}

# --- Touch the remaining slides so they get re-saved too ----------------
for ($i = 2; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $title = $slide.Shapes.Item(1)
    $title.Name = $title.Name
}
